$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.405.32"
$ws.Range("E2").Value = "  +6.26%  "

$ws.Range("D3").Value = "3.300.40"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'404.96"
$ws.Range("E5").Value = "  +3.17%  "

$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("D7").Value = "3.292.85"
$ws.Range("E7").Value = "  +2.50%  "

$ws.Range("E8").Value = "  -2.22%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").Value = "'0.613"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("E11").Value = "  +11.97%  "

$ws.Range("D12").Value = "'38.32"
$ws.Range("E12").Value = "  -1.09%  "

$ws.Range("D13").Value = "'0.143"
$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("D14").Value = "3.971.97"
$ws.Range("E14").Value = "  +6.04%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.07"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.472.79"
$ws.Range("E16").Value = "  +7.65%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.75"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "60.468.08"
$ws.Range("E18").Value = "  +6.49%  "

$ws.Range("D19").Value = "'0.985"
$ws.Range("E19").Value = "  -4.11%  "

$ws.Range("D20").Value = "'10.29"
$ws.Range("E20").Value = "  -6.38%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'3.21"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "'0.0000109"
$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'12.09"
$ws.Range("E23").Value = "  -6.71%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'289.21"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").Value = "'72.83"
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").Value = "'3.15"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.48"
$ws.Range("E27").Value = "  +2.05%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'28.25"
$ws.Range("E28").Value = "  +2.09%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'7.35"
$ws.Range("E30").Value = "  -3.67%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.162"
$ws.Range("E31").Value = "  -3.69%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'0.994"
$ws.Range("E32").Value = "  -0.59%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.108"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +12.31%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'37.42"
$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").Value = "'0.0471"
$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("D38").Value = "'51.85"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'27.60"
$ws.Range("E40").Value = "  +27.19%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("E42").Value = "  -7.92%  "

$ws.Range("D43").Value = "'136.92"
$ws.Range("E43").Value = "  +2.12%  "

$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").Value = "'1.84"
$ws.Range("E45").Value = "  -2.27%  "

$ws.Range("D46").Value = "'0.272"
$ws.Range("E46").Value = "  -2.86%  "

$ws.Range("E47").Value = "  -5.70%  "

$ws.Range("D48").Value = "'15.90"
$ws.Range("E48").Value = "  -5.38%  "

$ws.Range("D49").Value = "'2.22"
$ws.Range("E49").Value = "  +6.06%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.945.63"
$ws.Range("E50").Value = "  +10.82%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.126.51"
$ws.Range("E51").Value = "  -0.74%  "
